$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.998.18"
$ws.Range("E2").Value = "  +2.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.615.58"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.74"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.34"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.35"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.372"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.080.91"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.995.15"
$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.74"
$ws.Range("E16").Value = "  +4.49%  "

$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.630.65"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +3.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("E20").Value = "  +9.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.24"
$ws.Range("E21").Value = "  +4.18%  "

$ws.Range("E22").Value = "  +13.45%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.515"
$ws.Range("E24").Value = "  +12.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.26"
$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +1.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.76"
$ws.Range("E28").Value = "  +5.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0800"
$ws.Range("E29").Value = "  +1.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  +11.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").Value = "  +2.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.24"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("E34").Value = "  +2.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.957"
$ws.Range("E36").Value = "  +8.18%  "

$ws.Range("E37").Value = "  +4.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  +6.81%  "

$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.858"
$ws.Range("E40").Value = "  -2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "297.96"
$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.46"
$ws.Range("E43").Value = "  +11.02%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.608"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0553"
$ws.Range("E47").Value = "  +2.24%  "

$ws.Range("E48").Value = "  +3.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.73"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.69"
$ws.Range("E50").Value = "  +5.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.045.89"
$ws.Range("E51").Value = "  +4.62%  "
